$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3692.2222
$ws.Range("I74").Value = 3466.6667
$ws.Range("J74").Value = 3805
$ws.Range("K74").Value = 3466.6667
$ws.Range("L74").Value = 3805
$ws.Range("M74").Value = -2530.6667
$ws.Range("N74").Value = -5677
# Row 76
$ws.Range("H76").Value = 44141.793
$ws.Range("I76").Value = 44141.793
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 44141.793
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -43826.793
$ws.Range("N76").ClearContents()
# Row 77
$ws.Range("H77").Value = 3692.2222
$ws.Range("I77").Value = 3466.6667
$ws.Range("J77").Value = 3805
$ws.Range("K77").Value = 17333.3335
$ws.Range("L77").Value = 19025
$ws.Range("M77").Value = -12653.3335
$ws.Range("N77").Value = -28385
# Row 79
$ws.Range("H79").Value = 44141.793
$ws.Range("I79").Value = 44141.793
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 44141.793
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -43049.793
$ws.Range("N79").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2467.3
$ws.Range("I61").Value = 1486.1818
$ws.Range("J61").Value = 3666.4443
$ws.Range("K61").Value = 1486.1818
$ws.Range("L61").Value = 3666.4443
$ws.Range("M61").Value = -1274.1818
$ws.Range("N61").Value = -4090.4443
# Row 63
$ws.Range("H63").Value = 715535.7
$ws.Range("I63").Value = 770423.0600000001
$ws.Range("K63").Value = 770423.0600000001
$ws.Range("M63").Value = -769737.0600000001
# Row 66
$ws.Range("H66").Value = 715535.7
$ws.Range("I66").Value = 770423.0600000001
$ws.Range("K66").Value = 3852115.3
$ws.Range("M66").Value = -3848683.3
# Row 74
$ws.Range("H74").Value = 743.5217
$ws.Range("I74").Value = 647.55
$ws.Range("J74").Value = 1383.3334
$ws.Range("K74").Value = 647.55
$ws.Range("L74").Value = 1383.3334
$ws.Range("M74").Value = 226.45
$ws.Range("N74").Value = -3131.3334
# Row 77
$ws.Range("H77").Value = 743.5217
$ws.Range("I77").Value = 647.55
$ws.Range("J77").Value = 1383.3334
$ws.Range("K77").Value = 3237.75
$ws.Range("L77").Value = 6916.666999999999
$ws.Range("M77").Value = 1130.25
$ws.Range("N77").Value = -15652.667
# Row 132
$ws.Range("H132").Value = 5245.2563
$ws.Range("I132").Value = 6390.1665
$ws.Range("J132").Value = 3413.4
$ws.Range("K132").Value = 19170.4995
$ws.Range("L132").Value = 10240.2
$ws.Range("M132").Value = -16640.4995
$ws.Range("N132").Value = -15300.2
# Row 136
$ws.Range("H136").Value = 2467.3
$ws.Range("I136").Value = 1486.1818
$ws.Range("J136").Value = 3666.4443
$ws.Range("K136").Value = 4458.5454
$ws.Range("L136").Value = 10999.3329
$ws.Range("M136").Value = -1908.5454
$ws.Range("N136").Value = -16099.3329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 469.35294
$ws.Range("I22").Value = 467.4375
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 467.4375
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -294.4375
$ws.Range("N22").Value = -846
# Row 105
$ws.Range("H105").Value = 2609.6956
$ws.Range("I105").Value = 2006.5834
$ws.Range("J105").Value = 3267.6365
$ws.Range("K105").Value = 2006.5834
$ws.Range("L105").Value = 3267.6365
$ws.Range("M105").Value = -259.5834
$ws.Range("N105").Value = -6761.636500000001
# Row 134
$ws.Range("H134").Value = 21138.725
$ws.Range("I134").Value = 34618.7
$ws.Range("J134").Value = 1881.619
$ws.Range("K134").Value = 103856.1
$ws.Range("L134").Value = 5644.857
$ws.Range("M134").Value = -101321.1
$ws.Range("N134").Value = -10714.857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 83335310
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 83335310
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 83335310
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -83336558
# Row 65
$ws.Range("H65").Value = 83335310
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 83335310
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 416676550
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -416682790
# Row 132
$ws.Range("H132").Value = 2780255.8
$ws.Range("I132").Value = 2090.6086
$ws.Range("J132").Value = 5684701
$ws.Range("K132").Value = 6271.825800000001
$ws.Range("L132").Value = 17054103
$ws.Range("M132").Value = -3741.825800000001
$ws.Range("N132").Value = -17059163

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 28494362
$ws.Range("I70").Value = 38721540
$ws.Range("J70").Value = 4372
$ws.Range("K70").Value = 38721540
$ws.Range("L70").Value = 4372
$ws.Range("M70").Value = -38721270
$ws.Range("N70").Value = -4912
# Row 73
$ws.Range("H73").Value = 28494362
$ws.Range("I73").Value = 38721540
$ws.Range("J73").Value = 4372
$ws.Range("K73").Value = 38721540
$ws.Range("L73").Value = 4372
$ws.Range("M73").Value = -38720604
$ws.Range("N73").Value = -6244
# Row 80
$ws.Range("H80").Value = 4998.625
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 5897.8
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 5897.8
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -7893.8
# Row 83
$ws.Range("H83").Value = 4998.625
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 5897.8
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 29489
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -39473
# Row 132
$ws.Range("H132").Value = 94946.37
$ws.Range("I132").Value = 170068.67
$ws.Range("J132").Value = 4799.6
$ws.Range("K132").Value = 510206.01
$ws.Range("L132").Value = 14398.8
$ws.Range("M132").Value = -507676.01
$ws.Range("N132").Value = -19458.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 27779674
$ws.Range("I7").Value = 1157.1428
$ws.Range("J7").Value = 45456910
$ws.Range("K7").Value = 1157.1428
$ws.Range("L7").Value = 45456910
$ws.Range("M7").Value = -1045.1428
$ws.Range("N7").Value = -45457134
# Row 40
$ws.Range("H40").Value = 1432.0834
$ws.Range("I40").Value = 1432.0834
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1432.0834
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1296.0834
$ws.Range("N40").ClearContents()
# Row 126
$ws.Range("H126").Value = 27779674
$ws.Range("I126").Value = 1157.1428
$ws.Range("J126").Value = 45456910
$ws.Range("K126").Value = 3471.4284
$ws.Range("L126").Value = 136370730
$ws.Range("M126").Value = -1001.4284
$ws.Range("N126").Value = -136375670
# Row 132
$ws.Range("H132").Value = 7221.647
$ws.Range("I132").Value = 10125.714
$ws.Range("J132").Value = 2530.4614
$ws.Range("K132").Value = 30377.142
$ws.Range("L132").Value = 7591.3842
$ws.Range("M132").Value = -27847.142
$ws.Range("N132").Value = -12651.3842
